$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.021562707921193
$ws.Range("D2").Value = 1.026583254673848
$ws.Range("E2").Value = 1.022429520663156
$ws.Range("I2").Value = 1.030584396658172
$ws.Range("J2").Value = 1.026753473200626
$ws.Range("K2").Value = 1.029405477563145
$ws.Range("L2").Value = 1.025263927067645
$ws.Range("N2").Value = 1.013005452593201

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.022252695899515
$ws.Range("D3").Value = 1.027068441004793
$ws.Range("E3").Value = 1.023008187642047
$ws.Range("I3").Value = 1.030682306505229
$ws.Range("J3").Value = 1.027082484605793
$ws.Range("K3").Value = 1.029699082267485
$ws.Range("L3").Value = 1.025649881504466
$ws.Range("N3").Value = 1.013114311811413

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022699851145742
$ws.Range("D4").Value = 1.02738292320972
$ws.Range("E4").Value = 1.023383617840391
$ws.Range("I4").Value = 1.030744762724205
$ws.Range("J4").Value = 1.027295358324707
$ws.Range("K4").Value = 1.02988887672867
$ws.Range("L4").Value = 1.025899912433397
$ws.Range("N4").Value = 1.013184734933066

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.022887997494496
$ws.Range("D5").Value = 1.027515257423985
$ws.Range("E5").Value = 1.023541684532759
$ws.Range("I5").Value = 1.030770803590168
$ws.Range("J5").Value = 1.027384844612716
$ws.Range("K5").Value = 1.029968620114246
$ws.Range("L5").Value = 1.026005093959688
$ws.Range("N5").Value = 1.013214336500442

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022919597547912
$ws.Range("D6").Value = 1.027537484245304
$ws.Range("E6").Value = 1.023568238393736
$ws.Range("I6").Value = 1.030775163294054
$ws.Range("J6").Value = 1.027399869373266
$ws.Range("K6").Value = 1.029982006631071
$ws.Range("L6").Value = 1.026022758374618
$ws.Range("N6").Value = 1.013219306467392

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022702364531513
$ws.Range("D7").Value = 1.027384690973805
$ws.Range("E7").Value = 1.023385729011271
$ws.Range("I7").Value = 1.030745111532128
$ws.Range("J7").Value = 1.027296554069456
$ws.Range("K7").Value = 1.029889942446934
$ws.Range("L7").Value = 1.025901317606772
$ws.Range("N7").Value = 1.013185130488192

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.021795749067096
$ws.Range("D8").Value = 1.026747113857775
$ws.Range("E8").Value = 1.022624876796495
$ws.Range("I8").Value = 1.030617671237569
$ws.Range("J8").Value = 1.026864667194587
$ws.Range("K8").Value = 1.029504740649284
$ws.Range("L8").Value = 1.025394300675658
$ws.Range("N8").Value = 1.013042245123921

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020203533027826
$ws.Range("D9").Value = 1.02562780097102
$ws.Range("E9").Value = 1.021291864297752
$ws.Range("I9").Value = 1.030386260577968
$ws.Range("J9").Value = 1.026103544359989
$ws.Range("K9").Value = 1.028824589633989
$ws.Range("L9").Value = 1.024503180429445
$ws.Range("N9").Value = 1.012790360326784

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.019145786259638
$ws.Range("D10").Value = 1.024884521743532
$ws.Range("E10").Value = 1.020408497575334
$ws.Range("I10").Value = 1.030227429653652
$ws.Range("J10").Value = 1.025596154529087
$ws.Range("K10").Value = 1.028370310444167
$ws.Range("L10").Value = 1.02391074358789
$ws.Range("N10").Value = 1.012622396037276

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.018688682554096
$ws.Range("D11").Value = 1.024563393813974
$ws.Range("E11").Value = 1.020027276099971
$ws.Range("I11").Value = 1.030157582716422
$ws.Range("J11").Value = 1.025376472399893
$ws.Range("K11").Value = 1.028173419096872
$ws.Range("L11").Value = 1.023654621078227
$ws.Range("N11").Value = 1.012549661812781

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.018519032487275
$ws.Range("D12").Value = 1.024444222384059
$ws.Range("E12").Value = 1.019885868315556
$ws.Range("I12").Value = 1.030131478144978
$ws.Range("J12").Value = 1.025294877237608
$ws.Range("K12").Value = 1.028100258405471
$ws.Range("L12").Value = 1.023559548440324
$ws.Range("N12").Value = 1.012522644876901

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.018555416694497
$ws.Range("D13").Value = 1.024469780059184
$ws.Range("E13").Value = 1.019916191912632
$ws.Range("I13").Value = 1.030137084906255
$ws.Range("J13").Value = 1.025312379447947
$ws.Range("K13").Value = 1.028115952798129
$ws.Range("L13").Value = 1.023579938982692
$ws.Range("N13").Value = 1.012528440102981

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.018674656386257
$ws.Range("D14").Value = 1.024553540814904
$ws.Range("E14").Value = 1.020015583300813
$ws.Range("I14").Value = 1.030155428170252
$ws.Range("J14").Value = 1.025369727617154
$ws.Range("K14").Value = 1.028167372143453
$ws.Range("L14").Value = 1.023646761056794
$ws.Range("N14").Value = 1.012547428585865

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.018748142328876
$ws.Range("D15").Value = 1.024605163190083
$ws.Range("E15").Value = 1.020076847500519
$ws.Range("I15").Value = 1.030166708840892
$ws.Range("J15").Value = 1.025405062363375
$ws.Range("K15").Value = 1.028199049831879
$ws.Range("L15").Value = 1.023687940693284
$ws.Range("N15").Value = 1.012559128003897

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.019176142041135
$ws.Range("D16").Value = 1.024905849228011
$ws.Range("E16").Value = 1.020433825191904
$ws.Range("I16").Value = 1.030232042641963
$ws.Range("J16").Value = 1.025610734671227
$ws.Range("K16").Value = 1.028383373685291
$ws.Range("L16").Value = 1.023927750296415
$ws.Range("N16").Value = 1.01262722311287

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019444859532937
$ws.Range("D17").Value = 1.025094655034138
$ws.Range("E17").Value = 1.020658092700691
$ws.Range("I17").Value = 1.030272738315587
$ws.Range("J17").Value = 1.025739754095035
$ws.Range("K17").Value = 1.028498946494378
$ws.Range("L17").Value = 1.024078286491362
$ws.Range("N17").Value = 1.012669936490093

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.019601685321998
$ws.Range("D18").Value = 1.025204851296131
$ws.Range("E18").Value = 1.020789027704283
$ws.Range("I18").Value = 1.030296371920404
$ws.Range("J18").Value = 1.025815010875804
$ws.Range("K18").Value = 1.028566340171646
$ws.Range("L18").Value = 1.024166130827249
$ws.Range("N18").Value = 1.012694850006367

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.019655173623466
$ws.Range("D19").Value = 1.025242437019039
$ws.Range("E19").Value = 1.020833694047198
$ws.Range("I19").Value = 1.030304412791504
$ws.Range("J19").Value = 1.025840671778453
$ws.Range("K19").Value = 1.028589316580909
$ws.Range("L19").Value = 1.02419609006284
$ws.Range("N19").Value = 1.012703344776527

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019416019630519
$ws.Range("D20").Value = 1.025074390802939
$ws.Range("E20").Value = 1.020634018123064
$ws.Range("I20").Value = 1.030268382754238
$ws.Range("J20").Value = 1.025725911321754
$ws.Range("K20").Value = 1.028486548476224
$ws.Range("L20").Value = 1.02406213132448
$ws.Range("N20").Value = 1.01266535379451

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.018639539424919
$ws.Range("D21").Value = 1.02452887231287
$ws.Range("E21").Value = 1.019986309617796
$ws.Range("I21").Value = 1.030150030955394
$ws.Range("J21").Value = 1.025352839873495
$ws.Range("K21").Value = 1.028152231147677
$ws.Range("L21").Value = 1.023627081878861
$ws.Range("N21").Value = 1.012541836952768

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.018152138517721
$ws.Range("D22").Value = 1.024186519371049
$ws.Range("E22").Value = 1.01958019798126
$ws.Range("I22").Value = 1.030074691360735
$ws.Range("J22").Value = 1.025118302133273
$ws.Range("K22").Value = 1.027941880061365
$ws.Range("L22").Value = 1.023353912373794
$ws.Range("N22").Value = 1.012464176055213

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.018410442095841
$ws.Range("D23").Value = 1.024367946148073
$ws.Range("E23").Value = 1.019795377744922
$ws.Range("I23").Value = 1.030114717931316
$ws.Range("J23").Value = 1.02524263199129
$ws.Range("K23").Value = 1.028053405147051
$ws.Range("L23").Value = 1.023498689727195
$ws.Range("N23").Value = 1.012505345496633

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019429050871212
$ws.Range("D24").Value = 1.025083547123754
$ws.Range("E24").Value = 1.020644896006948
$ws.Range("I24").Value = 1.030270351164811
$ws.Range("J24").Value = 1.025732166269365
$ws.Range("K24").Value = 1.02849215066273
$ws.Range("L24").Value = 1.02406943102835
$ws.Range("N24").Value = 1.012667424518916

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020614510963147
$ws.Range("D25").Value = 1.025916662575403
$ws.Range("E25").Value = 1.021635553635517
$ws.Range("I25").Value = 1.03044689228811
$ws.Range("J25").Value = 1.026300314315779
$ws.Range("K25").Value = 1.029000579799313
$ws.Range("L25").Value = 1.024733273397821
$ws.Range("N25").Value = 1.012855487861789
